$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly refresh re-ordered the price observations (rows 2-11) for
# Damasco @ Vega Monumental Concepción. Columns A,B,C,E,F,G,H,I,J are
# identical for every row, so only D (Fecha), K (Variedad), L (Calidad),
# M (Volumen), N (Precio mínimo), O (Precio máximo), P (Precio promedio
# ponderado), Q (Unidad de comercialización), R (Origen), S (Precio $/Kg)
# and T (Kg / unidad) need to be rewritten per row.

$rows = @{
    2  = @{ D = 44875; K = "Castle Brite"; L = "Primera"; M = 50;  N = 31000; O = 32000; P = 31400; Q = "`$/bandeja 10 kilos";      R = "Provincia de Limarí";    S = 3140; T = 10 }
    3  = @{ D = 44559; K = "Modesto";      L = "Primera"; M = 100; N = 19000; O = 20000; P = 19500; Q = "`$/caja 18 kilos";         R = "Región de O'Higgins";    S = 1083; T = 18 }
    4  = @{ D = 44559; K = "Modesto";      L = "Segunda"; M = 50;  N = 18000; O = 18000; P = 18000; Q = "`$/caja 18 kilos";         R = "Región de O'Higgins";    S = 1000; T = 18 }
    5  = @{ D = 44545; K = "Castle Brite"; L = "Primera"; M = 100; N = 18000; O = 19000; P = 18500; Q = "`$/caja 15 kilos";         R = "Región de O'Higgins";    S = 1233; T = 15 }
    6  = @{ D = 44545; K = "Castle Brite"; L = "Segunda"; M = 50;  N = 17000; O = 17000; P = 17000; Q = "`$/caja 15 kilos";         R = "Región de O'Higgins";    S = 1133; T = 15 }
    7  = @{ D = 44579; K = "Modesto";      L = "Primera"; M = 180; N = 13000; O = 14000; P = 13444; Q = "`$/caja 18 kilos";         R = "Región Metropolitana";   S = 747;  T = 18 }
    8  = @{ D = 44159; K = "Castle Brite"; L = "Primera"; M = 100; N = 14000; O = 15000; P = 14500; Q = "`$/caja 15 kilos";         R = "Región Metropolitana";   S = 967;  T = 15 }
    9  = @{ D = 44189; K = "Dina";         L = "Primera"; M = 200; N = 15000; O = 16000; P = 15500; Q = "`$/caja 15 kilos granel";  R = "Región de O'Higgins";    S = 1033; T = 15 }
    10 = @{ D = 44189; K = "Dina";         L = "Segunda"; M = 100; N = 14000; O = 14000; P = 14000; Q = "`$/caja 15 kilos granel";  R = "Región de O'Higgins";    S = 933;  T = 15 }
    11 = @{ D = 44187; K = "Dina";         L = "Primera"; M = 100; N = 15000; O = 16000; P = 15500; Q = "`$/caja 18 kilos";         R = "Región Metropolitana";   S = 861;  T = 18 }
}

foreach ($r in $rows.Keys) {
    $data = $rows[$r]
    $ws.Cells.Item($r, 4).Value  = $data.D   # D - Fecha
    $ws.Cells.Item($r, 11).Value = $data.K   # K - Variedad
    $ws.Cells.Item($r, 12).Value = $data.L   # L - Calidad
    $ws.Cells.Item($r, 13).Value = $data.M   # M - Volumen
    $ws.Cells.Item($r, 14).Value = $data.N   # N - Precio minimo
    $ws.Cells.Item($r, 15).Value = $data.O   # O - Precio maximo
    $ws.Cells.Item($r, 16).Value = $data.P   # P - Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $data.Q   # Q - Unidad de comercializacion
    $ws.Cells.Item($r, 18).Value = $data.R   # R - Origen
    $ws.Cells.Item($r, 19).Value = $data.S   # S - Precio $/Kg
    $ws.Cells.Item($r, 20).Value = $data.T   # T - Kg / unidad
}
